# Auto-generated edit script applying numeric updates to Jenova_Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across all 8 sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 226.85715
$ws.Range("I2").Value = 217.8
$ws.Range("J2").Value = 249.5
$ws.Range("K2").Value = 217.8
$ws.Range("L2").Value = 249.5
$ws.Range("M2").Value = -104.8
$ws.Range("N2").Value = -475.5
$ws.Range("H17").Value = 42679.145
$ws.Range("J17").Value = 42679.145
$ws.Range("L17").Value = 128037.435
$ws.Range("N17").Value = -128373.435
$ws.Range("H40").Value = 4443.3335
$ws.Range("I40").Value = 3690.7693
$ws.Range("J40").Value = 6400
$ws.Range("K40").Value = 3690.7693
$ws.Range("L40").Value = 6400
$ws.Range("M40").Value = -3515.7693
$ws.Range("N40").Value = -6750
$ws.Range("H47").Value = 9487
$ws.Range("I47").Value = 4000
$ws.Range("J47").Value = 14974
$ws.Range("K47").Value = 4000
$ws.Range("L47").Value = 14974
$ws.Range("M47").Value = -3028
$ws.Range("N47").Value = -16918
$ws.Range("H62").Value = 12504356
$ws.Range("I62").Value = 15627207
$ws.Range("K62").Value = 15627207
$ws.Range("M62").Value = -15626583
$ws.Range("H65").Value = 12504356
$ws.Range("I65").Value = 15627207
$ws.Range("K65").Value = 78136035
$ws.Range("M65").Value = -78132915
$ws.Range("H69").Value = 9000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 9000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 27000
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -28748
$ws.Range("H72").Value = 9000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 9000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 81000
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -89736
$ws.Range("H106").Value = 2468.5454
$ws.Range("I106").Value = 2600.8667
$ws.Range("K106").Value = 2600.8667
$ws.Range("M106").Value = -1969.8667
$ws.Range("H107").Value = 63436.062
$ws.Range("I107").Value = 67611.8
$ws.Range("K107").Value = 67611.8
$ws.Range("M107").Value = -65691.8
$ws.Range("H137").Value = 3677.7925
$ws.Range("I137").Value = 2206.2727
$ws.Range("J137").Value = 6105.8
$ws.Range("K137").Value = 6618.8181
$ws.Range("L137").Value = 18317.4
$ws.Range("M137").Value = -4068.8181
$ws.Range("N137").Value = -23417.4
$ws.Range("H138").Value = 5278.245
$ws.Range("I138").Value = 3802.3333
$ws.Range("J138").Value = 5860.8423
$ws.Range("K138").Value = 11406.9999
$ws.Range("L138").Value = 17582.5269
$ws.Range("M138").Value = -6266.999899999999
$ws.Range("N138").Value = -27862.5269

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3050.3135
$ws.Range("I32").Value = 2673.7097
$ws.Range("K32").Value = 2673.7097
$ws.Range("M32").Value = -2386.7097
$ws.Range("H61").Value = 4999.6665
$ws.Range("I61").Value = 3428.1428
$ws.Range("K61").Value = 3428.1428
$ws.Range("M61").Value = -3216.1428
$ws.Range("H74").Value = 3935
$ws.Range("I74").Value = 2223.818
$ws.Range("K74").Value = 2223.818
$ws.Range("M74").Value = -1349.818
$ws.Range("H77").Value = 3935
$ws.Range("I77").Value = 2223.818
$ws.Range("K77").Value = 11119.09
$ws.Range("M77").Value = -6751.09
$ws.Range("H122").Value = 3952.9656
$ws.Range("I122").Value = 2689.45
$ws.Range("J122").Value = 6760.778
$ws.Range("K122").Value = 8068.349999999999
$ws.Range("L122").Value = 20282.334
$ws.Range("M122").Value = -5618.349999999999
$ws.Range("N122").Value = -25182.334
$ws.Range("H136").Value = 4999.6665
$ws.Range("I136").Value = 3428.1428
$ws.Range("K136").Value = 10284.4284
$ws.Range("M136").Value = -7734.428400000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 956.125
$ws.Range("I11").Value = 115
$ws.Range("J11").Value = 1797.25
$ws.Range("K11").Value = 115
$ws.Range("L11").Value = 1797.25
$ws.Range("M11").Value = 25
$ws.Range("N11").Value = -2077.25
$ws.Range("H20").Value = 45456412
$ws.Range("I20").Value = 71430050
$ws.Range("K20").Value = 71430050
$ws.Range("M20").Value = -71429803
$ws.Range("H26").Value = 9303.444
$ws.Range("I26").Value = 9303.444
$ws.Range("K26").Value = 9303.444
$ws.Range("M26").Value = -9011.444
$ws.Range("H28").Value = 29000
$ws.Range("J28").Value = 29000
$ws.Range("L28").Value = 29000
$ws.Range("N28").Value = -29588
$ws.Range("H107").Value = 2385.3635
$ws.Range("I107").Value = 1223.95
$ws.Range("K107").Value = 1223.95
$ws.Range("M107").Value = 696.05
$ws.Range("H134").Value = 4649.4287
$ws.Range("I134").Value = 3049.1428
$ws.Range("K134").Value = 9147.428400000001
$ws.Range("M134").Value = -6612.428400000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7600
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 7600
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 7600
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -8190
$ws.Range("H34").Value = 7600
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 7600
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 7600
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -8004
$ws.Range("H44").Value = 5000
$ws.Range("I44").Value = 5000
$ws.Range("K44").Value = 5000
$ws.Range("M44").Value = -4558
$ws.Range("H50").Value = 74950
$ws.Range("J50").Value = 74950
$ws.Range("L50").Value = 74950
$ws.Range("N50").Value = -76200
$ws.Range("H58").Value = 388948.06
$ws.Range("I58").Value = 1429971.2
$ws.Range("J58").Value = 5413.2104
$ws.Range("K58").Value = 1429971.2
$ws.Range("L58").Value = 5413.2104
$ws.Range("M58").Value = -1429768.2
$ws.Range("N58").Value = -5819.2104
$ws.Range("H76").Value = 9565
$ws.Range("I76").Value = 9565
$ws.Range("K76").Value = 9565
$ws.Range("M76").Value = -9250
$ws.Range("H79").Value = 9565
$ws.Range("I79").Value = 9565
$ws.Range("K79").Value = 9565
$ws.Range("M79").Value = -8473
$ws.Range("H132").Value = 4111.1816
$ws.Range("I132").Value = 3013.2942
$ws.Range("K132").Value = 9039.882599999999
$ws.Range("M132").Value = -6509.882599999999
$ws.Range("H134").Value = 4394.5
$ws.Range("I134").Value = 3557.8333
$ws.Range("K134").Value = 10673.4999
$ws.Range("M134").Value = -8138.499899999999
$ws.Range("H136").Value = 388948.06
$ws.Range("I136").Value = 1429971.2
$ws.Range("J136").Value = 5413.2104
$ws.Range("K136").Value = 4289913.6
$ws.Range("L136").Value = 16239.6312
$ws.Range("M136").Value = -4287363.6
$ws.Range("N136").Value = -21339.6312

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2331.6667
$ws.Range("I3").Value = 1997.5
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 5992.5
$ws.Range("L3").Value = 9000
$ws.Range("M3").Value = -5880.5
$ws.Range("N3").Value = -9224
$ws.Range("H129").Value = 2575.8572
$ws.Range("I129").Value = 675
$ws.Range("J129").Value = 3336.2
$ws.Range("K129").Value = 2025
$ws.Range("L129").Value = 10008.6
$ws.Range("M129").Value = 2975
$ws.Range("N129").Value = -20008.6

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 50000
$ws.Range("J42").Value = 50000
$ws.Range("L42").Value = 50000
$ws.Range("N42").Value = -50970
$ws.Range("H115").Value = 50000
$ws.Range("J115").Value = 50000
$ws.Range("L115").Value = 50000
$ws.Range("N115").Value = -52350
$ws.Range("H122").Value = 5620.9575
$ws.Range("I122").Value = 5526.533
$ws.Range("J122").Value = 5787.5884
$ws.Range("K122").Value = 16579.599
$ws.Range("L122").Value = 17362.7652
$ws.Range("M122").Value = -14129.599
$ws.Range("N122").Value = -22262.7652

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H100").Value = 12212.375
$ws.Range("I100").Value = 3957
$ws.Range("K100").Value = 3957
$ws.Range("M100").Value = -3416
$ws.Range("H107").Value = 3665.6667
$ws.Range("I107").Value = 3665.6667
$ws.Range("K107").Value = 3665.6667
$ws.Range("M107").Value = -1745.6667
$ws.Range("H136").Value = 4710.2856
$ws.Range("I136").Value = 3423.5715
$ws.Range("K136").Value = 10270.7145
$ws.Range("M136").Value = -7720.7145

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 7799.2
$ws.Range("J69").Value = 7799.2
$ws.Range("L69").Value = 7799.2
$ws.Range("N69").Value = -9297.200000000001
$ws.Range("H72").Value = 7799.2
$ws.Range("J72").Value = 7799.2
$ws.Range("L72").Value = 23397.6
$ws.Range("N72").Value = -30885.6
$ws.Range("H81").Value = 11656.615
$ws.Range("I81").Value = 1649.5
$ws.Range("K81").Value = 3299
$ws.Range("M81").Value = -2238
$ws.Range("H84").Value = 11656.615
$ws.Range("I84").Value = 1649.5
$ws.Range("K84").Value = 16495
$ws.Range("M84").Value = -11191
$ws.Range("H122").Value = 33337638
$ws.Range("I122").Value = 62503016
$ws.Range("J122").Value = 5777.9287
$ws.Range("K122").Value = 187509048
$ws.Range("L122").Value = 17333.7861
$ws.Range("M122").Value = -187506598
$ws.Range("N122").Value = -22233.7861
$ws.Range("H136").Value = 372662.72
$ws.Range("I136").Value = 402059.72
$ws.Range("K136").Value = 1206179.16
$ws.Range("M136").Value = -1203629.16
